$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain decimal numbers (e.g. "135.93").
# Excel would auto-convert these to numeric values on assignment, but the
# source data stores them as literal text (inline strings). Pre-formatting
# these cells as Text ("@") before assigning the value keeps them as text,
# matching the original file's representation.
$textCells = @("D6", "D9", "D11", "D12", "D14", "D20", "D21", "D23", "D25", "D27", "D29", "D30", "D32", "D35", "D36", "D37", "D40", "D43", "D45", "D48", "D51")
$textRange = $ws.Range($textCells[0])
for ($i = 1; $i -lt $textCells.Count; $i++) {
    $textRange = $excel.Union($textRange, $ws.Range($textCells[$i]))
}
foreach ($area in $textRange.Areas) {
    $area.NumberFormat = "@"
}

# Apply all cell value updates from the diff
$ws.Range("D2").Value = "62.872.40"
$ws.Range("E2").Value = "  -2.06%  "
$ws.Range("D3").Value = "3.114.76"
$ws.Range("E3").Value = "  -0.36%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("E5").Value = "  -2.04%  "
$ws.Range("D6").Value = "135.93"
$ws.Range("E6").Value = "  -4.64%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "3.112.27"
$ws.Range("E8").Value = "  -0.32%  "
$ws.Range("D9").Value = "0.514"
$ws.Range("E9").Value = "  -1.63%  "
$ws.Range("E10").Value = "  -3.90%  "
$ws.Range("D11").Value = "5.22"
$ws.Range("E11").Value = "  -2.90%  "
$ws.Range("D12").Value = "0.453"
$ws.Range("E12").Value = "  -3.00%  "
$ws.Range("D14").Value = "33.87"
$ws.Range("E14").Value = "  -3.47%  "
$ws.Range("D15").Value = "3.627.27"
$ws.Range("E15").Value = "  -0.54%  "
$ws.Range("E16").Value = "  +1.58%  "
$ws.Range("D17").Value = "62.887.97"
$ws.Range("E17").Value = "  -1.86%  "
$ws.Range("D18").Value = "3.121.71"
$ws.Range("E18").Value = "  +0.31%  "
$ws.Range("E19").Value = "  -3.24%  "
$ws.Range("D20").Value = "468.38"
$ws.Range("E20").Value = "  -1.95%  "
$ws.Range("D21").Value = "14.03"
$ws.Range("E21").Value = "  -3.40%  "
$ws.Range("E22").Value = "  -2.50%  "
$ws.Range("D23").Value = "7.62"
$ws.Range("E23").Value = "  -0.58%  "
$ws.Range("E24").Value = "  +0.17%  "
$ws.Range("D25").Value = "12.87"
$ws.Range("E25").Value = "  -3.71%  "
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("D27").Value = "2.70"
$ws.Range("E27").Value = "  -1.55%  "
$ws.Range("E28").Value = "  -6.78%  "
$ws.Range("B29").Value = "NEARProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D29").Value = "6.81"
$ws.Range("E29").Value = "  -4.82%  "
$ws.Range("B30").Value = "ImmutableX"
$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D30").Value = "2.07"
$ws.Range("E30").Value = "  +1.66%  "
$ws.Range("E31").Value = "  +0.05%  "
$ws.Range("D32").Value = "26.47"
$ws.Range("E32").Value = "  -1.23%  "
$ws.Range("E33").Value = "  -5.68%  "
$ws.Range("E34").Value = "  -5.13%  "
$ws.Range("D35").Value = "1.06"
$ws.Range("E35").Value = "  -3.41%  "
$ws.Range("B36").Value = "OKB"
$ws.Range("C36").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D36").Value = "51.86"
$ws.Range("E36").Value = "  -0.74%  "
$ws.Range("B37").Value = "Filecoin"
$ws.Range("C37").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D37").Value = "5.72"
$ws.Range("E37").Value = "  -3.81%  "
$ws.Range("D38").Value = "0.0₃0687"
$ws.Range("E38").Value = "  -10.27%  "
$ws.Range("E39").Value = "  -1.88%  "
$ws.Range("D40").Value = "416.22"
$ws.Range("E40").Value = "  -6.07%  "
$ws.Range("E41").Value = "  -0.31%  "
$ws.Range("D42").Value = "2.893.64"
$ws.Range("E42").Value = "  +1.52%  "
$ws.Range("D43").Value = "2.65"
$ws.Range("E43").Value = "  -12.10%  "
$ws.Range("E44").Value = "  -6.01%  "
$ws.Range("D45").Value = "0.258"
$ws.Range("E45").Value = "  -0.34%  "
$ws.Range("E46").Value = "  +0.01%  "
$ws.Range("E47").Value = "  -6.09%  "
$ws.Range("D48").Value = "25.22"
$ws.Range("E48").Value = "  -2.84%  "
$ws.Range("E49").Value = "  -0.54%  "
$ws.Range("E50").Value = "  -7.86%  "
$ws.Range("D51").Value = "120.33"
$ws.Range("E51").Value = "  +0.30%  "
